$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New wallet-topup rows captured since the last export (data updated till
# 17 Jan 2021 8AM). Two additional orders land in rows 19 and 20, right
# after the previous last data row (18). Copy the formatting from row 18
# first so the new rows inherit the same number formats / borders, then
# fill in the values and (re)apply the derived-column formulas so the
# shared-formula results recompute correctly.
# ---------------------------------------------------------------------------

$ws.Range("A18:E18").Copy() | Out-Null
$ws.Range("A19:E20").PasteSpecial(-4122) | Out-Null

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 26025571
$ws.Range("C19").Value = 52795
$ws.Range("D19").Value = 49999.89
$ws.Range("E19").Value = 44210
$ws.Range("F19").Formula = "=IF(B19="""","""",C19-D19)"
$ws.Range("G19").Formula = "=IF(B19="""","""",F19/D19*100)"
$ws.Range("H19").Formula = "=IF(B19="""","""",D19*1.04)"
$ws.Range("I19").Formula = "=IF(B19="""","""",C19-H19)"

# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 26040898
$ws.Range("C20").Value = 369564
$ws.Range("D20").Value = 350000.23
$ws.Range("E20").Value = 44212
$ws.Range("F20").Formula = "=IF(B20="""","""",C20-D20)"
$ws.Range("G20").Formula = "=IF(B20="""","""",F20/D20*100)"
$ws.Range("H20").Formula = "=IF(B20="""","""",D20*1.04)"
$ws.Range("I20").Formula = "=IF(B20="""","""",C20-H20)"

$wb.Application.Calculate()

# The author's last click before saving landed on E21.
$ws.Range("E21").Select() | Out-Null
